# Refresh the crypto price / 1h-volume columns (D, E) on the single sheet
# with the latest scraped values (GitHub Actions cron update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value into a cell while preventing Excel from
# auto-converting numeric-looking text (e.g. "1.00", "19.60") into a real
# number. We briefly force a Text number format, assign the value, then
# restore the cell style back to the default ("Normal") so no residual
# formatting is left behind on the cell.
function Set-TextValue {
    param($addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '44.957.46'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').Value = '2.270.16'
$ws.Range('E3').Value = '  +1.86%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  -0.55%  '
Set-TextValue 'D5' '301.78'
$ws.Range('E5').Value = '  -1.13%  '
Set-TextValue 'D6' '94.23'
$ws.Range('E6').Value = '  +0.33%  '
Set-TextValue 'D7' '0.565'
$ws.Range('E7').Value = '  -1.01%  '
Set-TextValue 'D8' '1.00'
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('E9').Value = '  -0.52%  '
Set-TextValue 'D10' '34.13'
$ws.Range('E10').Value = '  -1.43%  '
Set-TextValue 'D11' '0.0791'
$ws.Range('E11').Value = '  -0.65%  '
Set-TextValue 'D12' '7.22'
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('D14').Value = '2.620.82'
$ws.Range('E14').Value = '  +1.98%  '
$ws.Range('D15').Value = '2.271.90'
$ws.Range('E15').Value = '  +2.05%  '
Set-TextValue 'D16' '13.63'
$ws.Range('E16').Value = '  +1.16%  '
Set-TextValue 'D17' '0.801'
$ws.Range('E17').Value = '  -3.57%  '
$ws.Range('D18').Value = '44.918.46'
$ws.Range('E18').Value = '  +1.40%  '
Set-TextValue 'D19' '13.05'
$ws.Range('E19').Value = '  +9.80%  '
$ws.Range('D20').Value = '0.0₃0918'
$ws.Range('E20').Value = '  -2.25%  '
Set-TextValue 'D21' '6.05'
$ws.Range('E21').Value = '  -2.42%  '
Set-TextValue 'D22' '65.73'
$ws.Range('E22').Value = '  +1.06%  '
Set-TextValue 'D23' '238.72'
$ws.Range('E23').Value = '  +0.45%  '
$ws.Range('E24').Value = '  -1.81%  '
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('E26').Value = '  -3.09%  '
Set-TextValue 'D27' '41.07'
$ws.Range('E27').Value = '  +9.28%  '
$ws.Range('E28').Value = '  -0.32%  '
Set-TextValue 'D29' '9.61'
$ws.Range('E29').Value = '  -1.33%  '
Set-TextValue 'D30' '19.60'
$ws.Range('E30').Value = '  -0.99%  '
Set-TextValue 'D31' '152.05'
$ws.Range('E31').Value = '  +1.07%  '
Set-TextValue 'D32' '5.55'
$ws.Range('E32').Value = '  -6.55%  '
Set-TextValue 'D33' '0.0791'
$ws.Range('E33').Value = '  +0.53%  '
Set-TextValue 'D34' '2.55'
$ws.Range('E34').Value = '  -2.53%  '
Set-TextValue 'D35' '2.93'
$ws.Range('E35').Value = '  -3.41%  '
$ws.Range('E36').Value = '  -1.02%  '
$ws.Range('E37').Value = '  -2.88%  '
Set-TextValue 'D38' '1.77'
$ws.Range('E38').Value = '  -3.02%  '
$ws.Range('E39').Value = '  +6.07%  '
Set-TextValue 'D40' '0.0310'
$ws.Range('E40').Value = '  +3.79%  '
Set-TextValue 'D41' '3.22'
$ws.Range('E41').Value = '  -3.59%  '
Set-TextValue 'D42' '13.60'
$ws.Range('E42').Value = '  -10.31%  '
$ws.Range('E43').Value = '  -0.76%  '
Set-TextValue 'D44' '1.92'
$ws.Range('E44').Value = '  +12.15%  '
$ws.Range('D45').Value = '1.744.63'
$ws.Range('E45').Value = '  -4.49%  '
Set-TextValue 'D46' '0.194'
$ws.Range('E46').Value = '  +3.36%  '
Set-TextValue 'D47' '76.63'
$ws.Range('E47').Value = '  -3.76%  '
Set-TextValue 'D48' '69.26'
$ws.Range('E48').Value = '  +0.77%  '
Set-TextValue 'D49' '95.61'
$ws.Range('E49').Value = '  -2.35%  '
Set-TextValue 'D50' '53.57'
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('E51').Value = '  -2.99%  '
